$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '52.266.98'
$ws.Range("E2").Value = '  +1.69%  '

$ws.Range("D3").Value = '2.797.86'
$ws.Range("E3").Value = '  +1.61%  '

$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").Value = '''345.98'
$ws.Range("E5").Value = '  +4.43%  '

$ws.Range("D6").Value = '''116.83'
$ws.Range("E6").Value = '  +1.05%  '

$ws.Range("E7").Value = '  +4.04%  '

$ws.Range("E8").Value = '  -0.12%  '

$ws.Range("D9").Value = '''0.585'
$ws.Range("E9").Value = '  +2.44%  '

$ws.Range("D10").Value = '''43.45'
$ws.Range("E10").Value = '  +4.44%  '

$ws.Range("D11").Value = '''0.0858'
$ws.Range("E11").Value = '  +3.38%  '

$ws.Range("D12").Value = '''20.15'
$ws.Range("E12").Value = '  -0.76%  '

$ws.Range("E13").Value = '  +2.01%  '

$ws.Range("D14").Value = '''7.81'
$ws.Range("E14").Value = '  +1.49%  '

$ws.Range("D15").Value = '3.235.37'
$ws.Range("E15").Value = '  +1.72%  '

$ws.Range("D16").Value = '2.782.53'
$ws.Range("E16").Value = '  +1.53%  '

$ws.Range("D17").Value = '''0.893'
$ws.Range("E17").Value = '  +0.60%  '

$ws.Range("D18").Value = '52.142.60'
$ws.Range("E18").Value = '  +1.54%  '

$ws.Range("D19").Value = '''3.22'
$ws.Range("E19").Value = '  +6.18%  '

$ws.Range("D20").Value = '''7.13'
$ws.Range("E20").Value = '  +3.96%  '

$ws.Range("D21").Value = '''13.44'
$ws.Range("E21").Value = '  -1.26%  '

$ws.Range("E22").Value = '  +1.87%  '

$ws.Range("D23").Value = '''70.29'
$ws.Range("E23").Value = '  +0.10%  '

$ws.Range("D24").Value = '''270.16'
$ws.Range("E24").Value = '  -3.91%  '

$ws.Range("E25").Value = '  +6.27%  '

$ws.Range("D26").Value = '''26.68'
$ws.Range("E26").Value = '  -0.76%  '

$ws.Range("E27").Value = '  -0.11%  '

$ws.Range("D28").Value = '''10.27'
$ws.Range("E28").Value = '  -0.69%  '

$ws.Range("E29").Value = '  +0.46%  '

$ws.Range("E30").Value = '  -0.42%  '

$ws.Range("D31").Value = '''35.09'
$ws.Range("E31").Value = '  -1.66%  '

$ws.Range("D32").Value = '''50.28'
$ws.Range("E32").Value = '  +0.23%  '

$ws.Range("E33").Value = '  +1.34%  '

$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D34").Value = '''0.0824'
$ws.Range("E34").Value = '  -0.09%  '

$ws.Range("B35").Value = 'VeChain'
$ws.Range("C35").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D35").Value = '''0.0415'
$ws.Range("E35").Value = '  +17.44%  '

$ws.Range("D36").Value = '''2.13'
$ws.Range("E36").Value = '  +0.92%  '

$ws.Range("E37").Value = '  -0.06%  '

$ws.Range("D38").Value = '''18.95'
$ws.Range("E38").Value = '  -2.59%  '

$ws.Range("D39").Value = '''4.99'
$ws.Range("E39").Value = '  -0.90%  '

$ws.Range("E40").Value = '  +0.00%  '

$ws.Range("E41").Value = '  +20.84%  '

$ws.Range("D42").Value = '''23.68'
$ws.Range("E42").Value = '  -0.03%  '

$ws.Range("D43").Value = '''128.08'
$ws.Range("E43").Value = '  -0.91%  '

$ws.Range("E44").Value = '  +2.54%  '

$ws.Range("D45").Value = '''2.31'
$ws.Range("E45").Value = '  +0.27%  '

$ws.Range("D46").Value = '''3.35'
$ws.Range("E46").Value = '  -2.33%  '

$ws.Range("D47").Value = '2.077.64'
$ws.Range("E47").Value = '  -1.70%  '

$ws.Range("E48").Value = '  +5.44%  '

$ws.Range("D49").Value = '''0.983'
$ws.Range("E49").Value = '  +17.80%  '

$ws.Range("E50").Value = '  -0.42%  '

$ws.Range("D51").Value = '''8.96'
$ws.Range("E51").Value = '  -1.39%  '
